$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D header ("kaat") ---
$ws.Range("D1").Value = "kaat"

# --- Fill in missing Name (column B) values for rows 2-30 ---
$ws.Range("B2").Value = "jabir"
$ws.Range("B3").Value = "zakir"
$ws.Range("B4").Value = "amjad"
$ws.Range("B5").Value = "akram"
$ws.Range("B6").Value = "safi"
$ws.Range("B7").Value = "farookh"
$ws.Range("B8").Value = "gulfam"
$ws.Range("B9").Value = "aarif"
$ws.Range("B10").Value = "toheed"
$ws.Range("B11").Value = "intazar"
$ws.Range("B12").Value = "jabir"
$ws.Range("B13").Value = "akram"
$ws.Range("B14").Value = "amjad"
$ws.Range("B15").Value = "zakir"
$ws.Range("B16").Value = "farookh"
$ws.Range("B17").Value = "aarif"
$ws.Range("B18").Value = "gulfam"
$ws.Range("B19").Value = "toheed"
$ws.Range("B20").Value = "safi"
$ws.Range("B21").Value = "intazar"

# --- New row 22 (previously absent) ---
$ws.Range("B22").Value = "hasan"
$ws.Range("C22").Value = 360

# --- Row 23: fill Name + new kaat quantity ---
$ws.Range("B23").Value = "wazid"
$ws.Range("D23").Value = 20

# --- Fill in missing Name (column B) values for rows 24-30 ---
$ws.Range("B24").Value = "sehzad"
$ws.Range("B25").Value = "jabir"
$ws.Range("B26").Value = "akram"
$ws.Range("B27").Value = "amjad"
$ws.Range("B28").Value = "zakir"
$ws.Range("B29").Value = "farookh"
$ws.Range("B30").Value = "aarif"

# --- Fill previously-missing Quantity value for "hasan" rows ---
$ws.Range("C33").Value = 360
$ws.Range("D44").Value = 20
$ws.Range("C49").Value = 360
$ws.Range("C57").Value = 360
$ws.Range("D59").Value = 20

# --- Row 455 changes from "kayum" (no qty) to "gulfam" 450/50 ---
$ws.Range("B455").Value = "gulfam"
$ws.Range("C455").Value = 450
$ws.Range("D455").Value = 50

# --- New row 456 ---
$ws.Range("B456").Value = "hasan"
$ws.Range("C456").Value = 538

# --- Update the view: scroll position + active selection ---
$ws.Range("C457").Select()
$excel.ActiveWindow.ScrollRow = 433
$excel.ActiveWindow.ScrollColumn = 1
